# Insert a new data row at row 134 (pushing the existing rows 134-170 down
# to 135-171) and populate it with the new "Red Lady" price observation for
# Agrícola del Norte S.A. de Arica (Papa), per the weekly update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 134..170 down to 135..171, leaving a blank row 134 to fill in.
$ws.Rows.Item(134).Insert()

$ws.Range("A134").Value = 1
$ws.Range("B134").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C134").Value = "Arica y Parinacota"
$ws.Range("D134").Value = 44988
$ws.Range("E134").Value = 15
$ws.Range("F134").Value = 100114001
$ws.Range("G134").Value = "Papa"
$ws.Range("H134").Value = "Red Lady"
$ws.Range("I134").Value = "1a (cosecha)"
$ws.Range("J134").Value = 750
$ws.Range("K134").Value = 14000
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = 14600
$ws.Range("N134").Value = "$/saco 25 kilos"
$ws.Range("O134").Value = "Región Metropolitana"
$ws.Range("P134").Value = 584
$ws.Range("Q134").Value = 25
$ws.Range("R134").Value = "Hortaliza"
